$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text (e.g. "308.94", "26.961.83")
# formatted with thousands separators as literal dots, which Excel would
# otherwise reinterpret as a number. Mark the cells we are about to
# rewrite as Text first so the new values stay strings, matching the
# inline-string cell type used by the source workbook.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8:D37").NumberFormat = "@"
$ws.Range("D39:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.961.83'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.844.32'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '308.94'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  +2.06%  '
$ws.Range("D8").Value = '0.3677'
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("D9").Value = '0.07221'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").Value = '0.9312'
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").Value = '19.83'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '0.07725'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '1.863.11'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '5.379'
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '6.462'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '88.85'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '1.015'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = '0.000008659'
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").Value = '1.012'
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = '26.984.17'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").Value = '14.55'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").Value = '5.063'
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").Value = '10.64'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '1.955'
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").Value = '152.85'
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").Value = '18.24'
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = '2.005'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = '114.39'
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '4.965'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '0.08862'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = '3.307'
$ws.Range("E31").Value = '  +3.83%  '
$ws.Range("D32").Value = '1.176'
$ws.Range("E32").Value = '  -0.29%  '
$ws.Range("D33").Value = '0.7424'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '4.503'
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("D35").Value = '2.677'
$ws.Range("E35").Value = '  -5.86%  '
$ws.Range("D36").Value = '1.113'
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("D37").Value = '0.01963'
$ws.Range("E37").Value = '  +1.12%  '
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("D39").Value = '2.966'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").Value = '0.5243'
$ws.Range("E40").Value = '  +1.51%  '
$ws.Range("D41").Value = '7.010'
$ws.Range("E41").Value = '  +1.63%  '
$ws.Range("D42").Value = '0.1510'
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = '8.295'
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = '10.60'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").Value = '0.4729'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = '1.012'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").Value = '101.76'
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("D48").Value = '1.608'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '65.71'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").Value = '0.06065'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = '0.8918'
$ws.Range("E51").Value = '  +3.32%  '
